{"js": "// The document had a stray \"_GoBack\" bookmark (Word auto-inserts this at the\n// location of the most recent edit) sitting right after the sentence\n// \"There is no \"Maybe\" because our technology is not yet advanced enough \".\n// The author then made a further edit up in the \"Boolean data type\" paragraph,\n// so Word moved the \"_GoBack\" bookmark there instead - it now spans from\n// right after \" (usually \" up to right after \" of logic \" (just before \"and \").\n// Reproduce that: drop the old bookmark and re-create it at the new span.\n\nconst body = context.document.body;\n\n// 1) Remove the existing \"_GoBack\" bookmark (wherever it currently sits).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the new start boundary: right after the \" (usually \" text run.\nconst startHits = body.search(\"\\u00A0(usually \", { matchCase: true, matchWholeWord: false });\nstartHits.load(\"items\");\nawait context.sync();\nif (startHits.items.length === 0) {\n  throw new Error(\"Could not find the '(usually ' anchor text\");\n}\nconst startPoint = startHits.items[0].getRange(\"End\");\n\n// 3) Find the new end boundary: right after the \" of logic \" text (before \"and \").\nconst endHits = body.search(\"\\u00A0of logic \", { matchCase: true, matchWholeWord: false });\nendHits.load(\"items\");\nawait context.sync();\nif (endHits.items.length === 0) {\n  throw new Error(\"Could not find the 'of logic ' anchor text\");\n}\nconst endPoint = endHits.items[0].getRange(\"End\");\n\n// 4) Re-insert \"_GoBack\" spanning from the start boundary to the end boundary -\n//    this splits the two enclosing runs and drops a bookmarkStart/bookmarkEnd\n//    pair around \"denoted true and false), intended to represent the truth\n//    values of logic \".\nconst span = startPoint.expandTo(endPoint);\nspan.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document had a stray \"_GoBack\" bookmark (Word auto-drops this at the\n# location of the most recent edit) sitting right after the sentence\n# \"There is no \"Maybe\" because our technology is not yet advanced enough \".\n# The author made a further edit up in the \"Boolean data type\" paragraph, so\n# Word moved the \"_GoBack\" bookmark there instead - it now spans from right\n# after \" (usually \" up to right after \" of logic \" (just before \"and \").\n# Reproduce that: drop the old bookmark and re-create it at the new span.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the existing \"_GoBack\" bookmark (wherever it currently sits).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Find the new start boundary: right after the \" (usually \" text run.\n$startRng = $d.Content\n$startRng.Find.Text = [char]160 + \"(usually \"\n$startRng.Find.MatchCase = $true\n$startRng.Find.Execute() | Out-Null\n$startRng.Collapse(0)   # wdCollapseEnd\n$startPos = $startRng.Start\n\n# 3) Find the new end boundary: right after the \" of logic \" text (before \"and \").\n$endRng = $d.Content\n$endRng.Find.Text = [char]160 + \"of logic \"\n$endRng.Find.MatchCase = $true\n$endRng.Find.Execute() | Out-Null\n$endRng.Collapse(0)   # wdCollapseEnd\n$endPos = $endRng.Start\n\n# 4) Re-insert \"_GoBack\" spanning from the start boundary to the end boundary -\n#    this splits the two enclosing runs and drops a bookmarkStart/bookmarkEnd\n#    pair around \"denoted true and false), intended to represent the truth\n#    values of logic \".\n$bmRange = $d.Range($startPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
